$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 175, shifting existing rows 175:279 down to 176:280
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new record's data
$ws.Cells.Item(175, 1).Value = 10
$ws.Cells.Item(175, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(175, 3).Value = "La Araucanía"
$ws.Cells.Item(175, 4).Value = 44879
$ws.Cells.Item(175, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(175, 5).Value = 9
$ws.Cells.Item(175, 6).Value = 100112052
$ws.Cells.Item(175, 7).Value = "Albahaca"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 100
$ws.Cells.Item(175, 11).Value = 9000
$ws.Cells.Item(175, 12).Value = 9000
$ws.Cells.Item(175, 13).Value = 9000
$ws.Cells.Item(175, 14).Value = "$/paquete"
$ws.Cells.Item(175, 15).Value = "Región Metropolitana"
$ws.Cells.Item(175, 16).Value = 9000
$ws.Cells.Item(175, 17).Value = 1
$ws.Cells.Item(175, 18).Value = "Hortaliza"
